# Update "想去人数" (interested-count) figures across the three sheets
# that carry this data: 展览 (exhibitions), 演出 (performances), and
# 全部类型 (the combined/aggregated sheet).

$wb = $excel.ActiveWorkbook

# ---- 展览 (sheet "展览") ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 342
$ws.Range("F4").Value = 1248
$ws.Range("F5").Value = 277
$ws.Range("F9").Value = 144
$ws.Range("F10").Value = 3464
$ws.Range("F12").Value = 87
$ws.Range("F14").Value = 41
$ws.Range("F16").Value = 595
$ws.Range("F17").Value = 89
$ws.Range("F18").Value = 738
$ws.Range("F19").Value = 209
$ws.Range("F24").Value = 2606
$ws.Range("F25").Value = 5114
$ws.Range("F26").Value = 33
$ws.Range("F27").Value = 74
$ws.Range("F28").Value = 477
$ws.Range("F29").Value = 2225
$ws.Range("F30").Value = 283
$ws.Range("F31").Value = 2242
$ws.Range("F33").Value = 492
$ws.Range("F35").Value = 114
$ws.Range("F36").Value = 175
$ws.Range("F39").Value = 797

# ---- 演出 (sheet "演出") ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 72

# ---- 全部类型 (sheet "全部类型") ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 342
$ws.Range("F4").Value = 1249
$ws.Range("F5").Value = 277
$ws.Range("F9").Value = 144
$ws.Range("F10").Value = 3464
$ws.Range("F12").Value = 87
$ws.Range("F14").Value = 72
$ws.Range("F15").Value = 41
$ws.Range("F17").Value = 595
$ws.Range("F18").Value = 89
$ws.Range("F19").Value = 738
$ws.Range("F20").Value = 209
$ws.Range("F25").Value = 2606
$ws.Range("F26").Value = 5114
$ws.Range("F27").Value = 33
$ws.Range("F28").Value = 74
$ws.Range("F29").Value = 477
$ws.Range("F30").Value = 2226
$ws.Range("F31").Value = 283
$ws.Range("F32").Value = 2242
$ws.Range("F34").Value = 492
$ws.Range("F36").Value = 114
$ws.Range("F37").Value = 175
$ws.Range("F40").Value = 797
